# Actualización automática 2025-08-20 14:40:09
# Updates the advisor "ALMEIDA CUATIN JHONATHANN CARLOS" figures for
# "PUERTAS DE SEGURIDAD" in agosto: the monthly sale moves from 0 to
# 145.49, which ripples through the three report sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": raw sale amount for the group / count row ---
$wsGrupo = $wb.Sheets("VENTAS POR GRUPO")
$wsGrupo.Range("N2").Value = 145.49
$wsGrupo.Range("N34").Value = "1 de 32"

# --- Sheet "VENTA MENSUAL": agosto column totals ---
$wsMensual = $wb.Sheets("VENTA MENSUAL")
$wsMensual.Range("F2").Value = 275.07
$wsMensual.Range("F34").Value = 15237.12

# --- Sheet "CUMPLIMIENTO MENSUAL": per-product row + grand total row ---
$wsCumpl = $wb.Sheets("CUMPLIMIENTO MENSUAL")
$wsCumpl.Range("D17").Value = 145.49
$wsCumpl.Range("E17").Value = 196.51
$wsCumpl.Range("F17").Value = 0.4254093567251462

$wsCumpl.Range("D19").Value = 15357.2
$wsCumpl.Range("E19").Value = 16752.08107555788
$wsCumpl.Range("F19").Value = 0.478279160591053
